$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '63.163.10'
$ws.Cells.Item(2, 5).Value = '  +0.33%  '

$ws.Cells.Item(3, 4).Value = '2.477.41'
$ws.Cells.Item(3, 5).Value = '  +0.96%  '

$ws.Cells.Item(4, 5).Value = '  -0.04%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '577.62'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.65%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '146.75'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +0.44%  '

$ws.Cells.Item(7, 5).Value = '  +0.01%  '

$ws.Cells.Item(8, 5).Value = '  -0.26%  '

$ws.Cells.Item(9, 4).Value = '2.475.99'
$ws.Cells.Item(9, 5).Value = '  +0.92%  '

$ws.Cells.Item(10, 5).Value = '  +0.35%  '

$ws.Cells.Item(11, 5).Value = '  +1.64%  '

$ws.Cells.Item(12, 5).Value = '  +0.71%  '

$ws.Cells.Item(13, 5).Value = '  +0.24%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '28.57'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +4.66%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.0000180'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +1.27%  '

$ws.Cells.Item(16, 4).Value = '2.925.21'
$ws.Cells.Item(16, 5).Value = '  +0.94%  '

$ws.Cells.Item(17, 4).Value = '62.983.32'
$ws.Cells.Item(17, 5).Value = '  +0.19%  '

$ws.Cells.Item(18, 4).Value = '2.471.43'
$ws.Cells.Item(18, 5).Value = '  +0.99%  '

$ws.Cells.Item(19, 5).Value = '  +4.15%  '

$ws.Cells.Item(21, 2).Value = 'SuiNetwork'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '2.28'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +11.16%  '

$ws.Cells.Item(22, 2).Value = 'BitcoinCash'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '329.04'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +0.11%  '

$ws.Cells.Item(24, 5).Value = '  -0.02%  '

$ws.Cells.Item(25, 2).Value = 'Bittensor'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '677.81'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +6.69%  '

$ws.Cells.Item(26, 2).Value = 'Litecoin'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '66.24'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +0.85%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '9.73'
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +14.41%  '

$ws.Cells.Item(28, 4).Value = '0.0₃0998'
$ws.Cells.Item(28, 5).Value = '  +0.92%  '

$ws.Cells.Item(29, 5).Value = '  +3.36%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.997'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +372.92%  '

$ws.Cells.Item(31, 5).Value = '  +3.03%  '

$ws.Cells.Item(32, 5).Value = '  -1.68%  '

$ws.Cells.Item(33, 5).Value = '  +0.53%  '

$ws.Cells.Item(34, 5).Value = '  -3.51%  '

$ws.Cells.Item(35, 5).Value = '  +3.31%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.998'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -0.09%  '

$ws.Cells.Item(37, 5).Value = '  +0.68%  '

$ws.Cells.Item(38, 5).Value = '  +1.31%  '

$ws.Cells.Item(39, 5).Value = '  -0.76%  '

$ws.Cells.Item(40, 5).Value = '  +0.68%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '151.83'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -0.78%  '

$ws.Cells.Item(42, 5).Value = '  -1.79%  '

$ws.Cells.Item(43, 5).Value = '  -0.16%  '

$ws.Cells.Item(45, 5).Value = '  +10.72%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '154.53'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +6.35%  '

$ws.Cells.Item(47, 5).Value = '  +16.18%  '

$ws.Cells.Item(48, 5).Value = '  +0.18%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '20.63'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +0.63%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.607'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +0.75%  '

$ws.Cells.Item(51, 5).Value = '  -0.81%  '
